# GH ACTION Headlines refresh: advance "days since previous payout" (col G)
# by one and decrement "days until next payout" (col I) by one for every
# bond row, reflecting that one more day has elapsed since the last run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$colG = 7   # "Dni od poprzedniej wypłaty"
$colI = 9   # "Dni do następnej wypłaty"

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, $colG)
    $iCell = $ws.Cells.Item($r, $colI)

    $gVal = $gCell.Value()
    if ($null -ne $gVal) {
        $gCell.Value = $gVal + 1
    }

    $iVal = $iCell.Value()
    if ($null -ne $iVal) {
        $iCell.Value = $iVal - 1
    }
}
